$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is "Felipe" (felipefm.suporte@gmail.com). Mark his participation
# with "x" in the columns: Exploração Sites (C), Python (Eng. dados) (D),
# Dev Crawler (E), Python DataScience (H), Dashboards (J) and User Git (L).
$ws.Range("C8").Value = "x"
$ws.Range("D8").Value = "x"
$ws.Range("E8").Value = "x"
$ws.Range("H8").Value = "x"
$ws.Range("J8").Value = "x"
$ws.Range("L8").Value = "x"

# Leave the selection where the author left it when saving.
$ws.Range("K11").Select()
